# Updates cryptos list values (prices / volume %) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.749.22"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "1.660.08"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9999"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3818"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.75%  "

$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.11%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.243"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08213"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.526"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.434"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001237"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("D17").Value = "1.651.11"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.847"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.19%  "

$ws.Range("D24").Value = "23.751.53"
$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.513"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.058"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.215"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("D31").Value = "1.834.92"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.968"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.200"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.61%  "

$ws.Range("B34").Value = "FraxShare"
$ws.Range("C34").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.23%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.066"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02821"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2530"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.137"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08788"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07097"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7047"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6551"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("E46").Value = "  +1.84%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07959"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
